$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 108 (existing row 108 "subtotal" -> 109, row 109 "footer" -> 110)
$ws.Rows("108:108").Insert()

# New row 108 should look/format exactly like the data row above it (row 107):
# copy formats only so no new style entries are created.
$ws.Range("A107:N107").Copy()
$ws.Range("A108").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Recreate the same merged-cell layout as the data row above (B:G, H:K, L:M).
$ws.Range("B108:G108").Merge()
$ws.Range("H108:K108").Merge()
$ws.Range("L108:M108").Merge()

# Row 107 (item 104) becomes the newly added product.
$ws.Range("B107").Value2 = "مرطب شفاه لونا جوز هند ابيض"
$ws.Range("H107").Value2 = "0:0"
$ws.Range("L107").Value2 = 20

# Row 108 (item 105) keeps the data that used to belong to item 104 ("مناديل سولو سحب").
$ws.Range("A108").Value2 = 105
$ws.Range("B108").Value2 = "مناديل سولو سحب"
$ws.Range("H108").Value2 = "28:0"
$ws.Range("L108").Value2 = 45
$ws.Range("N108").Value2 = "1:0"

# Subtotal row (now row 109) increases by the new item's price.
$ws.Range("K109").Value2 = 4393.6400000000003

# Row heights: new data row matches the data-row height; the shifted footer
# row's height is recalculated slightly smaller by Excel after the edit.
$ws.Rows("108:108").RowHeight = 25.5
$ws.Rows("110:110").RowHeight = 16.5
